$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.695329432185133
$ws.Range("C2").Value = 0.2929308297521231
$ws.Range("D2").Value = 0.0702208808798801
$ws.Range("F2").Value = 1.789860563089604
$ws.Range("G2").Value = 0.002496826371860562
$ws.Range("I2").Value = 1.318122297531318
$ws.Range("L2").Value = 0.235109538552237
$ws.Range("M2").Value = 0.3455716035821652
$ws.Range("N2").Value = 1.804913798639703
$ws.Range("B3").Value = 1.584179911766284
$ws.Range("C3").Value = 0.2548042521924287
$ws.Range("D3").Value = 0.071006820089071
$ws.Range("F3").Value = 1.757118184208281
$ws.Range("G3").Value = 0.002502051419105781
$ws.Range("I3").Value = 1.312308957100889
$ws.Range("L3").Value = 0.2335045244996437
$ws.Range("M3").Value = 0.3300313452470505
$ws.Range("N3").Value = 1.824448230199536
$ws.Range("B4").Value = 1.516847234068166
$ws.Range("C4").Value = 0.2313942500149722
$ws.Range("D4").Value = 0.07151628605891958
$ws.Range("F4").Value = 1.73821281942196
$ws.Range("G4").Value = 0.002505428824946347
$ws.Range("I4").Value = 1.309517864352863
$ws.Range("L4").Value = 0.232630672873178
$ws.Range("M4").Value = 0.3206870864115317
$ws.Range("N4").Value = 1.837098399163125
$ws.Range("B5").Value = 1.489637413992398
$ws.Range("C5").Value = 0.2218537303085952
$ws.Range("D5").Value = 0.07173064438209309
$ws.Range("F5").Value = 1.730808643683105
$ws.Range("G5").Value = 0.002506847838651629
$ws.Range("I5").Value = 1.308575494075995
$ws.Range("L5").Value = 0.2323026408241446
$ws.Range("M5").Value = 0.3169287719147462
$ws.Range("N5").Value = 1.842418036781691
$ws.Range("B6").Value = 1.485133035985029
$ws.Range("C6").Value = 0.220269467835351
$ws.Range("D6").Value = 0.07176664543341182
$ws.Range("F6").Value = 1.729597257007526
$ws.Range("G6").Value = 0.002507086047712321
$ws.Range("I6").Value = 1.308430772984153
$ws.Range("L6").Value = 0.232249866666919
$ws.Range("M6").Value = 0.3163076957051132
$ws.Range("N6").Value = 1.843311289213108
$ws.Range("B7").Value = 1.516479347767756
$ws.Range("C7").Value = 0.2312655872212019
$ws.Range("D7").Value = 0.07151914967526096
$ws.Range("F7").Value = 1.73811175174589
$ws.Range("G7").Value = 0.002505447789534859
$ws.Range("I7").Value = 1.309504366505521
$ws.Range("L7").Value = 0.2326261352661447
$ws.Range("M7").Value = 0.320636200013702
$ws.Range("N7").Value = 1.837169475717268
$ws.Range("B8").Value = 1.656814857634913
$ws.Range("C8").Value = 0.2797843659869557
$ws.Range("D8").Value = 0.07048628053692063
$ws.Range("F8").Value = 1.778321313982801
$ws.Range("G8").Value = 0.002498592937026709
$ws.Range("I8").Value = 1.315955893055126
$ws.Range("L8").Value = 0.2345329578577804
$ws.Range("M8").Value = 0.3401722465751149
$ws.Range("N8").Value = 1.811512917323405
$ws.Range("B9").Value = 1.939320230395481
$ws.Range("C9").Value = 0.3749656821183862
$ws.Range("D9").Value = 0.06867495463134787
$ws.Range("F9").Value = 1.866756019772382
$ws.Range("G9").Value = 0.002486486534620597
$ws.Range("I9").Value = 1.33481676114026
$ws.Range("L9").Value = 0.2391585945273889
$ws.Range("M9").Value = 0.3800576899285062
$ws.Range("N9").Value = 1.766420301124633
$ws.Range("B10").Value = 2.151438187061558
$ws.Range("C10").Value = 0.4449755668678677
$ws.Range("D10").Value = 0.06747567386617703
$ws.Range("F10").Value = 1.937683368725885
$ws.Range("G10").Value = 0.002478397093989374
$ws.Range("I10").Value = 1.352509810998143
$ws.Range("L10").Value = 0.2430990028778695
$ws.Range("M10").Value = 0.4103367513348744
$ws.Range("N10").Value = 1.736492805074583
$ws.Range("B11").Value = 2.248951374553883
$ws.Range("C11").Value = 0.4768553222651803
$ws.Range("D11").Value = 0.06695885757825337
$ws.Range("F11").Value = 1.971267961037
$ws.Range("G11").Value = 0.002474889841032895
$ws.Range("I11").Value = 1.361403048452871
$ws.Range("L11").Value = 0.2450096892692386
$ws.Range("M11").Value = 0.4243267415156993
$ws.Range("N11").Value = 1.723577665742518
$ws.Range("B12").Value = 2.286025259523171
$ws.Range("C12").Value = 0.4889328695966242
$ws.Range("D12").Value = 0.06676730567591349
$ws.Range("F12").Value = 1.984177165702135
$ws.Range("G12").Value = 0.002473586414134599
$ws.Range("I12").Value = 1.36489300206631
$ws.Range("L12").Value = 0.2457502311420967
$ws.Range("M12").Value = 0.4296556474714848
$ws.Range("N12").Value = 1.718787967656453
$ws.Range("B13").Value = 2.278034145153924
$ws.Range("C13").Value = 0.4863315019410948
$ws.Range("D13").Value = 0.06680837464432265
$ws.Range("F13").Value = 1.981388395109832
$ws.Range("G13").Value = 0.002473866034364034
$ws.Range("I13").Value = 1.364135926476578
$ws.Range("L13").Value = 0.245589985616121
$ws.Range("M13").Value = 0.4285065820610896
$ws.Range("N13").Value = 1.719815016008667
$ws.Range("B14").Value = 2.251998497513171
$ws.Range("C14").Value = 0.4778488370644709
$ws.Range("D14").Value = 0.06694301508312961
$ws.Range("F14").Value = 1.972326162217115
$ws.Range("G14").Value = 0.002474782113160824
$ws.Range("I14").Value = 1.361687713630587
$ws.Range("L14").Value = 0.2450702731820087
$ws.Range("M14").Value = 0.4247645276086516
$ws.Range("N14").Value = 1.723181587818608
$ws.Range("B15").Value = 2.236070188042959
$ws.Range("C15").Value = 0.4726536835438537
$ws.Range("D15").Value = 0.06702602803489732
$ws.Range("F15").Value = 1.966800266641997
$ws.Range("G15").Value = 0.002475346449411927
$ws.Range("I15").Value = 1.360204061484851
$ws.Range("L15").Value = 0.2447541494525893
$ws.Range("M15").Value = 0.4224764790976323
$ws.Range("N15").Value = 1.725256875415255
$ws.Range("B16").Value = 2.145086085080209
$ws.Range("C16").Value = 0.4428928590261876
$ws.Range("D16").Value = 0.06751002927356531
$ws.Range("F16").Value = 1.935515232770598
$ws.Range("G16").Value = 0.002478629764266369
$ws.Range("I16").Value = 1.351945682656449
$ws.Range("L16").Value = 0.2429765140932858
$ws.Range("M16").Value = 0.4094268300423423
$ws.Range("N16").Value = 1.737350945653787
$ws.Range("B17").Value = 2.089532239462812
$ws.Range("C17").Value = 0.4246442865301674
$ws.Range("D17").Value = 0.06781432555604283
$ws.Range("F17").Value = 1.916662022559024
$ws.Range("G17").Value = 0.002480688099985312
$ws.Range("I17").Value = 1.347096337641801
$ws.Range("L17").Value = 0.2419162686720853
$ws.Range("M17").Value = 0.4014766972211703
$ws.Range("N17").Value = 1.744949632702756
$ws.Range("B18").Value = 2.05767504286996
$ws.Range("C18").Value = 0.4141511458567493
$ws.Range("D18").Value = 0.06799205276452724
$ws.Range("F18").Value = 1.905942237195518
$ws.Range("G18").Value = 0.002481888261309404
$ws.Range("I18").Value = 1.344386568395421
$ws.Range("L18").Value = 0.241317565877651
$ws.Range("M18").Value = 0.396924302407605
$ws.Range("N18").Value = 1.749385964124613
$ws.Range("B19").Value = 2.046905181536829
$ws.Range("C19").Value = 0.4105988325236467
$ws.Range("D19").Value = 0.06805269190008545
$ws.Range("F19").Value = 1.902333961219554
$ws.Range("G19").Value = 0.002482297412376263
$ws.Range("I19").Value = 1.343482706911686
$ws.Range("L19").Value = 0.2411167652919204
$ws.Range("M19").Value = 0.3953864235882136
$ws.Range("N19").Value = 1.750899311299733
$ws.Range("B20").Value = 2.095436107170656
$ws.Range("C20").Value = 0.4265865666625359
$ws.Range("D20").Value = 0.06778165268103287
$ws.Range("F20").Value = 1.918656124501069
$ws.Range("G20").Value = 0.00248046730468754
$ws.Range("I20").Value = 1.347604330404067
$ws.Range("L20").Value = 0.2420279823977936
$ws.Range("M20").Value = 0.402320899181241
$ws.Range("N20").Value = 1.744133929453092
$ws.Range("B21").Value = 2.259641779074173
$ws.Range("C21").Value = 0.4803402491334055
$ws.Range("D21").Value = 0.06690335497290256
$ws.Range("F21").Value = 1.974982749994467
$ws.Range("G21").Value = 0.002474512369319948
$ws.Range("I21").Value = 1.362403487855644
$ws.Range("L21").Value = 0.2452224636635378
$ws.Range("M21").Value = 0.4258628116813838
$ws.Range("N21").Value = 1.722189999621527
$ws.Range("B22").Value = 2.367821468725595
$ws.Range("C22").Value = 0.5155030956485689
$ws.Range("D22").Value = 0.06635356697506367
$ws.Range("F22").Value = 2.012912179344113
$ws.Range("G22").Value = 0.002470764346896501
$ws.Range("I22").Value = 1.372788712361967
$ws.Range("L22").Value = 0.2474093750628157
$ws.Range("M22").Value = 0.4414307472757386
$ws.Range("N22").Value = 1.708437211644167
$ws.Range("B23").Value = 2.310004696185842
$ws.Range("C23").Value = 0.496732858474445
$ws.Range("D23").Value = 0.06664477467952779
$ws.Range("F23").Value = 1.992565773799384
$ws.Range("G23").Value = 0.002472751617541855
$ws.Range("I23").Value = 1.367180395616245
$ws.Range("L23").Value = 0.2462331039866967
$ws.Range("M23").Value = 0.4331051474185585
$ws.Range("N23").Value = 1.715723298812655
$ws.Range("B24").Value = 2.09276671496832
$ws.Range("C24").Value = 0.4257084674593443
$ws.Range("D24").Value = 0.06779641542088122
$ws.Range("F24").Value = 1.91775421993529
$ws.Range("G24").Value = 0.002480567074037677
$ws.Range("I24").Value = 1.347374423413683
$ws.Range("L24").Value = 0.2419774428395343
$ws.Range("M24").Value = 0.4019391786655433
$ws.Range("N24").Value = 1.744502497916905
$ws.Range("B25").Value = 1.862101824502361
$ws.Range("C25").Value = 0.3492067708165223
$ws.Range("D25").Value = 0.06914194472155621
$ws.Range("F25").Value = 1.841794590392453
$ws.Range("G25").Value = 0.002489619578668367
$ws.Range("I25").Value = 1.32904444796921
$ws.Range("L25").Value = 0.2378121845652146
$ws.Range("M25").Value = 0.369097452960311
$ws.Range("N25").Value = 1.778058201409365
